$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Update the "Alternative Flow" bullet that talks about activity 6 -> 11.
#    The sentence (originally split across 3 runs with identical formatting)
#    reads:
#      "In activity 6 of Normal flow, if user do not select confirm print
#       user id and password system will provide main menu page to user."
#    and becomes:
#      "In activity 11 of Normal flow, if user do not select confirm print
#       user id and password system will provide main menu page to user."
#    Only the "6" -> "11" portion of the text actually changes; match on the
#    unique leading fragment (this exact phrase only occurs once in the
#    document - the sibling bullet without "print" is left untouched).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "In activity 6 of Normal flow, if user do not select confirm print",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In activity 11 of Normal flow, if user do not select confirm print",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 2) The "_GoBack" bookmark (Word's automatic "last edit" marker) moves from
#    the end of the Assumption bullet ("...add new cashier to the system.")
#    to sit right in the middle of the sentence we just edited, between
#    "use" and "r do not select..." - i.e. right after the word "use" in
#    "if user do not select".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "In activity 11 of Normal flow, if use",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0
) | Out-Null
$newBookmarkPos = $rng.End

# Remove the old bookmark (best effort) before re-adding it at the new spot.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$newRange = $d.Range($newBookmarkPos, $newBookmarkPos)
$d.Bookmarks.Add("_GoBack", $newRange)
